$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the "_generated" tag to the generated OpenScenario result folder names
$ws.Range("B1:Y1").Value = "KTH_pedestrian_autoware_light/OpenScenario/Results/Experiment_A1_generated/OpenScenario"
$ws.Range("Z1:AW1").Value = "KTH_pedestrian_autoware_light/OpenScenario/Results/Experiment_A2_generated/OpenScenario"

# Widen the columns so the longer text still fits.
# (82.109375 is the on-disk "characters" width we're targeting; the host
# quantizes ColumnWidth to its internal pixel grid, so 81.25 is the input
# that lands on the closest representable value, ~82.1667.)
$ws.Range("B1:AW1").EntireColumn.ColumnWidth = 81.25
